# Update "jung-wien-alle-without-gaps_pivot" sheet:
#  - Insert a new "Hermann Bahr – Paul Goldmann" pair column (shifts the existing
#    Hermann Bahr–Richard Beer-Hofmann / Hugo von Hofmannsthal–Richard Beer-Hofmann /
#    Paul Goldmann–Richard Beer-Hofmann columns one place to the right).
#  - Update the "Felix Salten – Hermann Bahr" counts (Bahr <-> Salten, and vice versa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before the current column K ("Hermann Bahr – Richard
#    Beer-Hofmann"), pushing K/L/M to L/M/N.
$ws.Columns("K:K").Insert()

# 2. Header for the newly inserted column.
$ws.Range("K1").Value = "Hermann Bahr – Paul Goldmann"

# 3. Fill in the data for the new column (all zero except 1907, which is 4),
#    and refresh the "Felix Salten – Hermann Bahr" column with the corrected
#    counts.
$gValues = @(0,0,0,0,0,0,0,6,4,6,5,8,6,8,2,7,18,14,17,4,4,8,6,2,2,2,0,0,0,2,2,0,0,2,0,0,0,3,0,4,0,0,6,0)
$kValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,4,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt 44; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $gValues[$i]
    $ws.Range("K$row").Value = $kValues[$i]
}
